$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Cases sheet: drop the "tags" / "description" columns (G, H) that
#    are being replaced by the new standalone "Calls" sheet.
# ---------------------------------------------------------------------
$cases = $wb.Worksheets.Item("Cases")
$cases.Range("G1:H3").Delete(-4159)

# Whole-row selection on row 1 (matches the other non-active sheets in
# this workbook, e.g. Deals) and drops Cases as the active/visible tab.
$null = $cases.Rows.Item(1).Select()

# ---------------------------------------------------------------------
# 2. Add the new "Calls" sheet right after "Cases" - it becomes the
#    active tab.
# ---------------------------------------------------------------------
$new = $wb.Worksheets.Add($null, $cases)
$new.Name = "Calls"

# Column B ("flag") filled top to bottom first.
$new.Range("B1").Value = "flag"
$new.Range("B2").Value = "Busy"
$new.Range("B3").Value = "Call Back"

# Rest of the header row.
$new.Range("C1").Value = "deal"
$new.Range("D1").Value = "task"
$new.Range("E1").Value = "case"
$new.Range("F1").Value = "notes"

# Remaining data, column by column.
$new.Range("D2").Value = "cccc"
$new.Range("D3").Value = "dddd"
$new.Range("E2").Value = "eeee"
$new.Range("E3").Value = "ffff"
$new.Range("F2").Value = "gggg"
$new.Range("F3").Value = "hhhh"

# Column A ("contact") and the rest of column C ("deal").
$new.Range("A1").Value = "contact"
$new.Range("A2").Value = "zzzx"
$new.Range("A3").Value = "aaaa"
$new.Range("C2").Value = "aaaa"
$new.Range("C3").Value = "bbbb"

# Header row gets the same highlight fill used by the other sheets.
$new.Range("A1:F1").Interior.Color = $cases.Range("A1").Interior.Color
